# "calculated rank and box office $ correlations"
#
# Adds a new "Correlation" worksheet (after "Weekly Data") that pulls the
# weekly NYT bestseller Rank and weekly Box Office Gross figures from the
# "Weekly Data" sheet (columns C and D, rows 50-72) and computes the
# correlation between them with CORREL().

$wb = $excel.ActiveWorkbook

# --- capture the previous selection on "Weekly Data" ------------------
# The author had been looking at the bottom of the Weekly Data sheet
# (around the Box Office Gross totals) right before adding the new sheet.
$wsWeekly = $wb.Worksheets.Item("Weekly Data")
$wsWeekly.Activate()
$wsWeekly.Range("D68").Select()

# --- add the new sheet, right after "Weekly Data" ----------------------
$wsCorr = $wb.Worksheets.Add($null, $wb.Sheets("Weekly Data"))
$wsCorr.Name = "Correlation"

# --- headers -------------------------------------------------------------
$wsCorr.Range("A1").Value = "Rank"
$wsCorr.Range("B1").Value = "Box Office Gross"
$wsCorr.Range("D1").Value = "Correlation:"

# --- data: Rank (A) and Box Office Gross (B), one row per week ---------
$rank = @(1,1,1,1,1,2,2,1,1,1,1,1,1,1,2,2,2,2,4,3,4,5,5)
$gross = @(40397446,57919989,39712036,24044930,17372256,15198257,11159987,6296543,5866786,2451332,2054096,1263487,961058,1101545,595502,595033,567017,322524,222276,115750,107216,48296,22585)

for ($i = 0; $i -lt $rank.Length; $i++) {
    $row = $i + 2
    $wsCorr.Cells.Item($row, 1).Value = $rank[$i]
    $wsCorr.Cells.Item($row, 2).Value = $gross[$i]
}

# --- correlation formula -------------------------------------------------
$wsCorr.Range("D2").Formula = "=CORREL(A2:A24,B2:B24)"

# --- column B is wide enough to show the full gross figures ------------
$wsCorr.Columns.Item(2).ColumnWidth = 14

# --- leave the new sheet active, with D3 selected -----------------------
$wsCorr.Range("D3").Select()
